# Update the class schedule ("visualize groups") for the two
# "matematyka stosowana_1_1" / "matematyka stosowana_1_2" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "matematyka stosowana_1_1" ---
$ws1 = $wb.Worksheets.Item("matematyka stosowana_1_1")

# Row 2
$ws1.Range("C2").Value = $null
$ws1.Range("D2").Value = "Analiza_matematyczna_I_lecture_1"
$ws1.Range("E2").Value = "Analiza_matematyczna_I_practicals_2_grp_1"
$ws1.Range("F2").Value = $null

# Row 3
$ws1.Range("B3").Value = "Analiza_matematyczna_I_practicals_1_grp_1"
$ws1.Range("D3").Value = "Technologie_informatyczne_I_laboratories_1_grp_1"
$ws1.Range("F3").Value = $null

# Row 4
$ws1.Range("C4").Value = "Algebra_liniowa_z_geometrią_analityczną_I_lecture_1"
$ws1.Range("D4").Value = "Wstęp_do_logiki_i_teorii_mnogości_lecture_1"
$ws1.Range("F4").Value = "Algebra_liniowa_z_geometrią_analityczną_I_practicals_1_grp_1"

# Row 5
$ws1.Range("B5").Value = "Wstęp_do_logiki_i_teorii_mnogości_practicals_1_grp_1"
$ws1.Range("D5").Value = $null

# Row 6
$ws1.Range("B6").Value = "Wstęp_do_obliczeń_symbolicznych_laboratories_1_grp_1"
$ws1.Range("C6").Value = "Analiza_matematyczna_I_lecture_2"
$ws1.Range("E6").Value = $null

# Row 7
$ws1.Range("C7").Value = $null
$ws1.Range("F7").Value = $null

# --- Sheet "matematyka stosowana_1_2" ---
$ws2 = $wb.Worksheets.Item("matematyka stosowana_1_2")

# Row 2
$ws2.Range("C2").Value = $null
$ws2.Range("D2").Value = "Analiza_matematyczna_I_lecture_1"
$ws2.Range("E2").Value = "Analiza_matematyczna_I_practicals_2_grp_2"
$ws2.Range("F2").Value = $null

# Row 4
$ws2.Range("B4").Value = $null
$ws2.Range("C4").Value = "Algebra_liniowa_z_geometrią_analityczną_I_lecture_1"
$ws2.Range("D4").Value = "Wstęp_do_logiki_i_teorii_mnogości_lecture_1"

# Row 5
$ws2.Range("C5").Value = "Analiza_matematyczna_I_practicals_1_grp_2"
$ws2.Range("F5").Value = $null

# Row 6
$ws2.Range("B6").Value = "Wstęp_do_obliczeń_symbolicznych_laboratories_1_grp_2"
$ws2.Range("C6").Value = "Analiza_matematyczna_I_lecture_2"
$ws2.Range("D6").Value = $null
$ws2.Range("E6").Value = $null
$ws2.Range("F6").Value = "Technologie_informatyczne_I_laboratories_1_grp_2"

# Row 7
$ws2.Range("E7").Value = "Wstęp_do_logiki_i_teorii_mnogości_practicals_1_grp_2"
$ws2.Range("F7").Value = "Algebra_liniowa_z_geometrią_analityczną_I_practicals_1_grp_2"
